# Clear the PASS/SKIP result markers left over from the last test run.
$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsRoundTrip = $wb.Worksheets.Item("RoundTrip")

# "Test Cases" sheet: D3 held "PASS" -> clear it out.
$wsTestCases.Range("D3").Value = $null

# "RoundTrip" sheet: I2 held "PASS", I3 held "SKIP" -> clear both.
$wsRoundTrip.Range("I2").Value = $null
$wsRoundTrip.Range("I3").Value = $null

# Move the active selection on "Test Cases" from A3 to D3, without
# disturbing which sheet/tab is actually active (RoundTrip stays active).
[void]$wsTestCases.Activate()
[void]$wsTestCases.Range("D3").Select()
[void]$wsRoundTrip.Activate()
